$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed) date column for rows 2-10 from 45185 to 45204
$ws.Range("C2:C10").Value = 45204
